# Se procesan de nuevo los datos con las nuevas dimensiones curadas
# Update the measure/dimension metadata rows (2-4) for columns E..I so that
# they describe measures (iaest-measure:*) instead of dimensions
# (sdmx-dimension:*/iaest-dimension:*), and clear out the now-unused
# mapping file references in row 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: qualified names -> switch from dimension refs to measure refs
$ws.Range("E2").Value = "iaest-measure:residencia-comarca-nombre"
$ws.Range("F2").Value = "iaest-measure:sexo"
$ws.Range("G2").Value = "iaest-measure:residencia-provincia-nombre"
$ws.Range("H2").Value = "iaest-measure:residencia-ccaa-nombre"
$ws.Range("I2").Value = "iaest-measure:relacion-lugar-de-residencia-y-nacimiento"

# Row 3: role -> from "dim" to "medida"
$ws.Range("E3").Value = "medida"
$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "medida"
$ws.Range("H3").Value = "medida"
$ws.Range("I3").Value = "medida"

# Row 4: type -> from "skos:Concept"/URI-* to "xsd:int"
$ws.Range("E4").Value = "xsd:int"
$ws.Range("F4").Value = "xsd:int"
$ws.Range("G4").Value = "xsd:int"
$ws.Range("H4").Value = "xsd:int"
$ws.Range("I4").Value = "xsd:int"

# Row 5: clear the now-obsolete mapping file references (fully remove the
# cells, not just their contents, so the row matches the curated layout)
$ws.Range("F5").Clear()
$ws.Range("H5").Clear()
$ws.Range("I5").Clear()
